$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-17 06:23:43"
$wsZh.Range("G2").Value = "2016-01-17 06:24:27"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-17 06:23:53"
$wsDe.Range("G2").Value = "2016-01-17 06:24:45"
